$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in new Person / Room Number values for the "Garbage & Area around" row (row 14)
$ws.Range("B14").Value = "Luke C. and Bogdan T."
$ws.Range("C14").Value = "U 09 & U 16"

# Update the active selection to match the saved workbook state
$ws.Range("G17").Select()
